$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$c = $ws.Range("A1717")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1717")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1717")
$c.NumberFormat = "@"
$c.Value = 'DUS'
$c.Style = "Normal"
$c = $ws.Range("D1717")
$c.NumberFormat = "@"
$c.Value = '4510.11'
$c.Style = "Normal"
$c = $ws.Range("E1717")
$c.NumberFormat = "@"
$c.Value = 'M1'
$c.Style = "Normal"
$c = $ws.Range("F1717")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1717")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1717")
$c.NumberFormat = "@"
$c.Value = '$ 50'
$c.Style = "Normal"
$c = $ws.Range("I1717")
$c.NumberFormat = "@"
$c.Value = '$ 25'
$c.Style = "Normal"

$c = $ws.Range("A1718")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1718")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1718")
$c.NumberFormat = "@"
$c.Value = '1st Speed 1 Yr School >35mphm4'
$c.Style = "Normal"
$c = $ws.Range("D1718")
$c.NumberFormat = "@"
$c.Value = '4511.21B1A'
$c.Style = "Normal"
$c = $ws.Range("E1718")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1718")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1718")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1718")
$c.NumberFormat = "@"
$c.Value = '$ 50'
$c.Style = "Normal"
$c = $ws.Range("I1718")
$c.NumberFormat = "@"
$c.Value = '$ 33'
$c.Style = "Normal"

$c = $ws.Range("A1719")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1719")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1719")
$c.NumberFormat = "@"
$c.Value = 'Reckless Operation 1st In 1 Yr'
$c.Style = "Normal"
$c = $ws.Range("D1719")
$c.NumberFormat = "@"
$c.Value = '4511.20'
$c.Style = "Normal"
$c = $ws.Range("E1719")
$c.NumberFormat = "@"
$c.Value = 'MM'
$c.Style = "Normal"
$c = $ws.Range("F1719")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1719")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1719")
$c.NumberFormat = "@"
$c.Value = '$ 22'
$c.Style = "Normal"
$c = $ws.Range("I1719")
$c.NumberFormat = "@"
$c.Value = '$ 11'
$c.Style = "Normal"

$c = $ws.Range("A1720")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1720")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1720")
$c.NumberFormat = "@"
$c.Value = 'DUS'
$c.Style = "Normal"
$c = $ws.Range("D1720")
$c.NumberFormat = "@"
$c.Value = '4510.11'
$c.Style = "Normal"
$c = $ws.Range("E1720")
$c.NumberFormat = "@"
$c.Value = 'M1'
$c.Style = "Normal"
$c = $ws.Range("F1720")
$c.NumberFormat = "@"
$c.Value = 'Dismissed'
$c.Style = "Normal"
$c = $ws.Range("H1720")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("I1720")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("J1720")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("K1720")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"

$c = $ws.Range("A1721")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1721")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1721")
$c.NumberFormat = "@"
$c.Value = '1st Speed 1 Yr School >35mphm4'
$c.Style = "Normal"
$c = $ws.Range("D1721")
$c.NumberFormat = "@"
$c.Value = '4511.21B1A'
$c.Style = "Normal"
$c = $ws.Range("E1721")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1721")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1721")
$c.NumberFormat = "@"
$c.Value = 'Guilty - Allied Offense'
$c.Style = "Normal"
$c = $ws.Range("H1721")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1721")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1721")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1721")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1722")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1722")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1722")
$c.NumberFormat = "@"
$c.Value = 'Reckless Operation 1st In 1 Yr - AMENDED to Disorderly Conduct - Persistent'
$c.Style = "Normal"
$c = $ws.Range("D1722")
$c.NumberFormat = "@"
$c.Value = '2917.11(A)(1)'
$c.Style = "Normal"
$c = $ws.Range("E1722")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1722")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1722")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1722")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1722")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1722")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1722")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1723")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1723")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1723")
$c.NumberFormat = "@"
$c.Value = 'DUS'
$c.Style = "Normal"
$c = $ws.Range("D1723")
$c.NumberFormat = "@"
$c.Value = '4510.11'
$c.Style = "Normal"
$c = $ws.Range("E1723")
$c.NumberFormat = "@"
$c.Value = 'M1'
$c.Style = "Normal"
$c = $ws.Range("F1723")
$c.NumberFormat = "@"
$c.Value = 'Dismissed'
$c.Style = "Normal"
$c = $ws.Range("H1723")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("I1723")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("J1723")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("K1723")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"

$c = $ws.Range("A1724")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1724")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1724")
$c.NumberFormat = "@"
$c.Value = '1st Speed 1 Yr School >35mphm4'
$c.Style = "Normal"
$c = $ws.Range("D1724")
$c.NumberFormat = "@"
$c.Value = '4511.21B1A'
$c.Style = "Normal"
$c = $ws.Range("E1724")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1724")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1724")
$c.NumberFormat = "@"
$c.Value = 'Guilty - Allied Offense'
$c.Style = "Normal"
$c = $ws.Range("H1724")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1724")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1724")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1724")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1725")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1725")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1725")
$c.NumberFormat = "@"
$c.Value = 'Reckless Operation 1st In 1 Yr - AMENDED to Disorderly Conduct - Persistent'
$c.Style = "Normal"
$c = $ws.Range("D1725")
$c.NumberFormat = "@"
$c.Value = '2917.11(A)(1)'
$c.Style = "Normal"
$c = $ws.Range("E1725")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1725")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1725")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1725")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1725")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1725")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1725")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1726")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1726")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1726")
$c.NumberFormat = "@"
$c.Value = 'Possession of Marijuana less than 100 grams'
$c.Style = "Normal"
$c = $ws.Range("D1726")
$c.NumberFormat = "@"
$c.Value = '2925.11(C)(3)(a)'
$c.Style = "Normal"
$c = $ws.Range("E1726")
$c.NumberFormat = "@"
$c.Value = 'MM'
$c.Style = "Normal"
$c = $ws.Range("F1726")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1726")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1726")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1726")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1726")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1726")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1727")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1727")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1727")
$c.NumberFormat = "@"
$c.Value = 'DUS'
$c.Style = "Normal"
$c = $ws.Range("D1727")
$c.NumberFormat = "@"
$c.Value = '4510.11'
$c.Style = "Normal"
$c = $ws.Range("E1727")
$c.NumberFormat = "@"
$c.Value = 'M1'
$c.Style = "Normal"
$c = $ws.Range("F1727")
$c.NumberFormat = "@"
$c.Value = 'Dismissed'
$c.Style = "Normal"
$c = $ws.Range("H1727")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("I1727")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("J1727")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("K1727")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"

$c = $ws.Range("A1728")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1728")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1728")
$c.NumberFormat = "@"
$c.Value = '1st Speed 1 Yr School >35mphm4'
$c.Style = "Normal"
$c = $ws.Range("D1728")
$c.NumberFormat = "@"
$c.Value = '4511.21B1A'
$c.Style = "Normal"
$c = $ws.Range("E1728")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1728")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1728")
$c.NumberFormat = "@"
$c.Value = 'Guilty - Allied Offense'
$c.Style = "Normal"
$c = $ws.Range("H1728")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1728")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1728")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1728")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1729")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1729")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1729")
$c.NumberFormat = "@"
$c.Value = 'Reckless Operation 1st In 1 Yr - AMENDED to Disorderly Conduct - Persistent'
$c.Style = "Normal"
$c = $ws.Range("D1729")
$c.NumberFormat = "@"
$c.Value = '2917.11(A)(1)'
$c.Style = "Normal"
$c = $ws.Range("E1729")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1729")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1729")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1729")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1729")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1729")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1729")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1730")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1730")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1730")
$c.NumberFormat = "@"
$c.Value = 'Possession of Marijuana less than 100 grams'
$c.Style = "Normal"
$c = $ws.Range("D1730")
$c.NumberFormat = "@"
$c.Value = '2925.11(C)(3)(a)'
$c.Style = "Normal"
$c = $ws.Range("E1730")
$c.NumberFormat = "@"
$c.Value = 'MM'
$c.Style = "Normal"
$c = $ws.Range("F1730")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1730")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1730")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1730")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1730")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1730")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1731")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1731")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1731")
$c.NumberFormat = "@"
$c.Value = 'DUS'
$c.Style = "Normal"
$c = $ws.Range("D1731")
$c.NumberFormat = "@"
$c.Value = '4510.11'
$c.Style = "Normal"
$c = $ws.Range("E1731")
$c.NumberFormat = "@"
$c.Value = 'M1'
$c.Style = "Normal"
$c = $ws.Range("F1731")
$c.NumberFormat = "@"
$c.Value = 'Dismissed'
$c.Style = "Normal"
$c = $ws.Range("H1731")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("I1731")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("J1731")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("K1731")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"

$c = $ws.Range("A1732")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1732")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1732")
$c.NumberFormat = "@"
$c.Value = '1st Speed 1 Yr School >35mphm4'
$c.Style = "Normal"
$c = $ws.Range("D1732")
$c.NumberFormat = "@"
$c.Value = '4511.21B1A'
$c.Style = "Normal"
$c = $ws.Range("E1732")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1732")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1732")
$c.NumberFormat = "@"
$c.Value = 'Guilty - Allied Offense'
$c.Style = "Normal"
$c = $ws.Range("H1732")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1732")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1732")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"
$c = $ws.Range("K1732")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1733")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1733")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1733")
$c.NumberFormat = "@"
$c.Value = 'Reckless Operation 1st In 1 Yr - AMENDED to Disorderly Conduct - Persistent'
$c.Style = "Normal"
$c = $ws.Range("D1733")
$c.NumberFormat = "@"
$c.Value = '2917.11(A)(1)'
$c.Style = "Normal"
$c = $ws.Range("E1733")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1733")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1733")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1733")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1733")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1733")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1733")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1734")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1734")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1734")
$c.NumberFormat = "@"
$c.Value = 'Possession of Marijuana less than 100 grams'
$c.Style = "Normal"
$c = $ws.Range("D1734")
$c.NumberFormat = "@"
$c.Value = '2925.11(C)(3)(a)'
$c.Style = "Normal"
$c = $ws.Range("E1734")
$c.NumberFormat = "@"
$c.Value = 'MM'
$c.Style = "Normal"
$c = $ws.Range("F1734")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1734")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1734")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1734")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1734")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1734")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1735")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1735")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1735")
$c.NumberFormat = "@"
$c.Value = 'DUS'
$c.Style = "Normal"
$c = $ws.Range("D1735")
$c.NumberFormat = "@"
$c.Value = '4510.11'
$c.Style = "Normal"
$c = $ws.Range("E1735")
$c.NumberFormat = "@"
$c.Value = 'M1'
$c.Style = "Normal"
$c = $ws.Range("F1735")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1735")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1735")
$c.NumberFormat = "@"
$c.Value = '$ 342'
$c.Style = "Normal"
$c = $ws.Range("I1735")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1735")
$c.NumberFormat = "@"
$c.Value = '12'
$c.Style = "Normal"
$c = $ws.Range("K1735")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1736")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1736")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1736")
$c.NumberFormat = "@"
$c.Value = '1st Speed 1 Yr School >35mphm4'
$c.Style = "Normal"
$c = $ws.Range("D1736")
$c.NumberFormat = "@"
$c.Value = '4511.21B1A'
$c.Style = "Normal"
$c = $ws.Range("E1736")
$c.NumberFormat = "@"
$c.Value = 'M4'
$c.Style = "Normal"
$c = $ws.Range("F1736")
$c.NumberFormat = "@"
$c.Value = 'No Contest'
$c.Style = "Normal"
$c = $ws.Range("G1736")
$c.NumberFormat = "@"
$c.Value = 'Guilty'
$c.Style = "Normal"
$c = $ws.Range("H1736")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("I1736")
$c.NumberFormat = "@"
$c.Value = '$ 0'
$c.Style = "Normal"
$c = $ws.Range("J1736")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"
$c = $ws.Range("K1736")
$c.NumberFormat = "@"
$c.Value = 'None'
$c.Style = "Normal"

$c = $ws.Range("A1737")
$c.NumberFormat = "@"
$c.Value = '21TRD09437'
$c.Style = "Normal"
$c = $ws.Range("B1737")
$c.NumberFormat = "@"
$c.Value = 'Hemmeter'
$c.Style = "Normal"
$c = $ws.Range("C1737")
$c.NumberFormat = "@"
$c.Value = 'Reckless Operation 1st In 1 Yr'
$c.Style = "Normal"
$c = $ws.Range("D1737")
$c.NumberFormat = "@"
$c.Value = '4511.20'
$c.Style = "Normal"
$c = $ws.Range("E1737")
$c.NumberFormat = "@"
$c.Value = 'MM'
$c.Style = "Normal"
$c = $ws.Range("F1737")
$c.NumberFormat = "@"
$c.Value = 'Dismissed'
$c.Style = "Normal"
$c = $ws.Range("H1737")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("I1737")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("J1737")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"
$c = $ws.Range("K1737")
$c.NumberFormat = "@"
$c.Value = ' '
$c.Style = "Normal"

$c = $ws.Range("A1738")
$c.NumberFormat = "@"
$c.Value = '21CRB00626'
$c.Style = "Normal"
$c = $ws.Range("B1738")
$c.NumberFormat = "@"
$c.Value = 'Rohrer'
$c.Style = "Normal"
$c = $ws.Range("C1738")
$c.NumberFormat = "@"
$c.Value = 'Criminal Mischief'
$c.Style = "Normal"
$c = $ws.Range("D1738")
$c.NumberFormat = "@"
$c.Value = '2909.07(A)(1)'
$c.Style = "Normal"
$c = $ws.Range("E1738")
$c.NumberFormat = "@"
$c.Value = 'M3'
$c.Style = "Normal"
$c = $ws.Range("F1738")
$c.NumberFormat = "@"
$c.Value = 'Not Guilty'
$c.Style = "Normal"

$c = $ws.Range("A1739")
$c.NumberFormat = "@"
$c.Value = '21CRB00626'
$c.Style = "Normal"
$c = $ws.Range("B1739")
$c.NumberFormat = "@"
$c.Value = 'Rohrer'
$c.Style = "Normal"
$c = $ws.Range("C1739")
$c.NumberFormat = "@"
$c.Value = 'Assault - M1'
$c.Style = "Normal"
$c = $ws.Range("D1739")
$c.NumberFormat = "@"
$c.Value = '2903.13(A)'
$c.Style = "Normal"
$c = $ws.Range("E1739")
$c.NumberFormat = "@"
$c.Value = 'M1'
$c.Style = "Normal"
$c = $ws.Range("F1739")
$c.NumberFormat = "@"
$c.Value = 'Not Guilty'
$c.Style = "Normal"

$c = $ws.Range("A1740")
$c.NumberFormat = "@"
$c.Value = '21CRB00626'
$c.Style = "Normal"
$c = $ws.Range("B1740")
$c.NumberFormat = "@"
$c.Value = 'Rohrer'
$c.Style = "Normal"
$c = $ws.Range("C1740")
$c.NumberFormat = "@"
$c.Value = 'Disorderly Conduct'
$c.Style = "Normal"
$c = $ws.Range("D1740")
$c.NumberFormat = "@"
$c.Value = '2917.11A1'
$c.Style = "Normal"
$c = $ws.Range("E1740")
$c.NumberFormat = "@"
$c.Value = 'MM'
$c.Style = "Normal"
$c = $ws.Range("F1740")
$c.NumberFormat = "@"
$c.Value = 'Not Guilty'
$c.Style = "Normal"

$c = $ws.Range("A1741")
$c.NumberFormat = "@"
$c.Value = '21CRB00626'
$c.Style = "Normal"
$c = $ws.Range("B1741")
$c.NumberFormat = "@"
$c.Value = 'Rohrer'
$c.Style = "Normal"
$c = $ws.Range("C1741")
$c.NumberFormat = "@"
$c.Value = 'Criminal Mischief - Victim is Family or Household Member'
$c.Style = "Normal"
$c = $ws.Range("D1741")
$c.NumberFormat = "@"
$c.Value = '2909.07**'
$c.Style = "Normal"
$c = $ws.Range("E1741")
$c.NumberFormat = "@"
$c.Value = 'M3'
$c.Style = "Normal"
$c = $ws.Range("F1741")
$c.NumberFormat = "@"
$c.Value = 'Not Guilty'
$c.Style = "Normal"

Write-Output ("Used range: " + $ws.UsedRange.Address())